# Update the single data row (row 2) on the active sheet to reflect the
# corrected customer transaction: name, email, per-flavor pizza counts and
# the recomputed order total.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "die"
$ws.Range("C2").Value = "nes"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 109.4
